$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 666757.6
$ws.Range("I33").Value = 909181.8
$ws.Range("K33").Value = 909181.8
$ws.Range("M33").Value = -908952.8
# Row 40
$ws.Range("H40").Value = 4169577
$ws.Range("J40").Value = 9094353
$ws.Range("L40").Value = 9094353
$ws.Range("N40").Value = -9094703
# Row 51
$ws.Range("H51").Value = 4929
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
# Row 55
$ws.Range("H55").Value = 136.86363
$ws.Range("I55").Value = 170.08333
$ws.Range("J55").Value = 97
$ws.Range("K55").Value = 170.08333
$ws.Range("L55").Value = 97
$ws.Range("M55").Value = 43.91667000000001
$ws.Range("N55").Value = -525
# Row 64
$ws.Range("H64").Value = 8427.571
$ws.Range("I64").Value = 8049.75
$ws.Range("J64").Value = 8931.333000000001
$ws.Range("K64").Value = 8049.75
$ws.Range("L64").Value = 8931.333000000001
$ws.Range("M64").Value = -7801.75
$ws.Range("N64").Value = -9427.333000000001
# Row 67
$ws.Range("H67").Value = 8427.571
$ws.Range("I67").Value = 8049.75
$ws.Range("J67").Value = 8931.333000000001
$ws.Range("K67").Value = 8049.75
$ws.Range("L67").Value = 8931.333000000001
$ws.Range("M67").Value = -7191.75
$ws.Range("N67").Value = -10647.333
# Row 100
$ws.Range("H100").Value = 3106.3635
$ws.Range("I100").Value = 1990
$ws.Range("K100").Value = 1990
$ws.Range("M100").Value = -1449

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 75
$ws.Range("H75").Value = 15000
$ws.Range("I75").Value = 15000
$ws.Range("K75").Value = 15000
$ws.Range("M75").Value = -14126
# Row 78
$ws.Range("H78").Value = 15000
$ws.Range("I78").Value = 15000
$ws.Range("K78").Value = 45000
$ws.Range("M78").Value = -40632
# Row 88
$ws.Range("H88").Value = 1094.3636
$ws.Range("I88").Value = 998.75
$ws.Range("J88").Value = 1149
$ws.Range("K88").Value = 998.75
$ws.Range("L88").Value = 1149
$ws.Range("M88").Value = -592.75
$ws.Range("N88").Value = -1961
# Row 91
$ws.Range("H91").Value = 1094.3636
$ws.Range("I91").Value = 998.75
$ws.Range("J91").Value = 1149
$ws.Range("K91").Value = 998.75
$ws.Range("L91").Value = 1149
$ws.Range("M91").Value = 405.25
$ws.Range("N91").Value = -3957

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 253.54546
$ws.Range("I22").Value = 254.33333
$ws.Range("K22").Value = 254.33333
$ws.Range("M22").Value = -81.33332999999999
# Row 86
$ws.Range("H86").Value = 2894.2292
$ws.Range("I86").Value = 2643.3513
$ws.Range("K86").Value = 2643.3513
$ws.Range("M86").Value = -1520.3513
# Row 89
$ws.Range("H89").Value = 2894.2292
$ws.Range("I89").Value = 2643.3513
$ws.Range("K89").Value = 13216.7565
$ws.Range("M89").Value = -7600.7565

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 4843.909
$ws.Range("I62").Value = 4838.6
$ws.Range("J62").Value = 4848.3335
$ws.Range("K62").Value = 4838.6
$ws.Range("L62").Value = 4848.3335
$ws.Range("M62").Value = -4214.6
$ws.Range("N62").Value = -6096.3335
# Row 65
$ws.Range("H65").Value = 4843.909
$ws.Range("I65").Value = 4838.6
$ws.Range("J65").Value = 4848.3335
$ws.Range("K65").Value = 24193
$ws.Range("L65").Value = 24241.6675
$ws.Range("M65").Value = -21073
$ws.Range("N65").Value = -30481.6675

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 4127.5
$ws.Range("I39").Value = 850
$ws.Range("J39").Value = 4783
$ws.Range("K39").Value = 2550
$ws.Range("L39").Value = 14349
$ws.Range("M39").Value = -2256
$ws.Range("N39").Value = -14937
# Row 94
$ws.Range("H94").Value = 5699
$ws.Range("I94").Value = 1400
$ws.Range("J94").Value = 9998
$ws.Range("K94").Value = 4200
$ws.Range("L94").Value = 29994
$ws.Range("M94").Value = -3524
$ws.Range("N94").Value = -31346
# Row 121
$ws.Range("H121").Value = 100404.5
$ws.Range("I121").Value = 125280.75
$ws.Range("J121").Value = 899.5
$ws.Range("K121").Value = 375842.25
$ws.Range("L121").Value = 2698.5
$ws.Range("M121").Value = -374532.25
$ws.Range("N121").Value = -5318.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Range("H27").Value = 4000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -834
# Row 31
$ws.Range("H31").Value = 9631.5
$ws.Range("I31").Value = 12265.5
$ws.Range("J31").Value = 6997.5
$ws.Range("K31").Value = 12265.5
$ws.Range("L31").Value = 6997.5
$ws.Range("M31").Value = -11973.5
$ws.Range("N31").Value = -7581.5
# Row 37
$ws.Range("H37").Value = 9631.5
$ws.Range("I37").Value = 12265.5
$ws.Range("J37").Value = 6997.5
$ws.Range("K37").Value = 12265.5
$ws.Range("L37").Value = 6997.5
$ws.Range("M37").Value = -11988.5
$ws.Range("N37").Value = -7551.5
# Row 70
$ws.Range("H70").Value = 22423.682
$ws.Range("I70").Value = 33865.19
$ws.Range("K70").Value = 33865.19
$ws.Range("M70").Value = -33595.19
# Row 73
$ws.Range("H73").Value = 22423.682
$ws.Range("I73").Value = 33865.19
$ws.Range("K73").Value = 33865.19
$ws.Range("M73").Value = -32929.19
# Row 102
$ws.Range("H102").Value = 1449.2153
$ws.Range("I102").Value = 1382.36
$ws.Range("K102").Value = 1382.36
$ws.Range("M102").Value = 239.6400000000001
# Row 126
$ws.Range("H126").Value = 2964.1072
$ws.Range("I126").Value = 2929.1304
$ws.Range("K126").Value = 8787.3912
$ws.Range("M126").Value = -6317.3912

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3865.0667
$ws.Range("I68").Value = 2335.182
$ws.Range("J68").Value = 8072.25
$ws.Range("K68").Value = 2335.182
$ws.Range("L68").Value = 8072.25
$ws.Range("M68").Value = -1586.182
$ws.Range("N68").Value = -9570.25
# Row 71
$ws.Range("H71").Value = 3865.0667
$ws.Range("I71").Value = 2335.182
$ws.Range("J71").Value = 8072.25
$ws.Range("K71").Value = 11675.91
$ws.Range("L71").Value = 40361.25
$ws.Range("M71").Value = -7931.91
$ws.Range("N71").Value = -47849.25

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1858.0322
$ws.Range("I132").Value = 1604.3478
$ws.Range("J132").Value = 2587.375
$ws.Range("K132").Value = 4813.0434
$ws.Range("L132").Value = 7762.125
$ws.Range("M132").Value = -2283.0434
$ws.Range("N132").Value = -12822.125
# Row 136
$ws.Range("H136").Value = 6090.7915
$ws.Range("I136").Value = 6098.864
$ws.Range("J136").Value = 6002
$ws.Range("K136").Value = 18296.592
$ws.Range("L136").Value = 18006
$ws.Range("M136").Value = -15746.592
$ws.Range("N136").Value = -23106
